$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18: new entry (clone the date cell format from A17, then set the value)
$ws.Range("A17").Copy($ws.Range("A18"))
$ws.Range("A18").Value = 42906
$ws.Range("B18").Value = 1.5
$ws.Range("C18").Value = "Inverzní logo, dodání kurzu, galerie, obrázky v seznamu služeb a jednotlivbých službách, ičo"

# Row 19: new entry (clone the date cell format from A17, then set the value)
$ws.Range("A17").Copy($ws.Range("A19"))
$ws.Range("A19").Value = 42908
$ws.Range("B19").Value = 0.3
$ws.Range("C19").Value = "Doplnění popisků galerie, singulár/plurál"

# Update selection to C23
$ws.Range("C23").Select()

$wb.Save()
